$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "My Users" -> "LR"
$ws.Name = "LR"

# Row 2 field edits
$ws.Range("C2").Value = "SS"
$ws.Range("D2").Value = "SS"
$ws.Range("E2").Value = "S"

# F2 must stay textual "13" (not be auto-converted to a number).
# Temporarily force text format so the numeric-looking string isn't
# coerced to a number, then drop back to the default "Normal" style so
# the cell keeps no explicit style override (matches the original).
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "13"
$ws.Range("F2").Style = "Normal"

# G2 becomes a plain number (was a shared-string "12")
$ws.Range("G2").Value = 1234567

# H2 numeric update
$ws.Range("H2").Value = 123

# Drop the now-unused trailing numeric cells on row 2
$ws.Range("I2:K2").ClearContents()

# Remove the second data row entirely
$ws.Rows(3).Delete()

Write-Host "done"
